# Adds a "Whatsearch"/"Button" search-bar pair of columns (D,E) before the
# existing "End" column (which slides from D to F), fills in the new
# "Mobil" value and a button xpath styled in a small monospace font, widens
# column D, moves the active selection to E2, and sets the print page setup
# to A4/portrait - matching the "searchbar; +whatsearch; +button click" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old D column ("End" header / empty data cell) moves out to F to make
# room for the two new columns.
$ws.Range("F1").Value = "End"

# New header row cells.
$ws.Range("D1").Value = "Whatsearch"
$ws.Range("E1").Value = "Button"

# New data row cells.
$ws.Range("D2").Value = "Mobil"
$ws.Range("E2").Value = '//*[@id="rootHead"]/form/button[2]'

# E2 (the button xpath) gets its own small monospace font.
$ws.Range("E2").Font.Name = "Consolas"
$ws.Range("E2").Font.Family = 3
$ws.Range("E2").Font.Size = 7
$ws.Range("E2").Font.Color = 2367776   # RGB(0x20, 0x21, 0x24) == FF202124

# Column D is widened so "Whatsearch"/"Mobil" fit comfortably.
$ws.Columns("D").ColumnWidth = 11.166666666666666   # renders as width="12"

# Selection moves to the newly added button-xpath cell.
$ws.Range("E2").Select() | Out-Null

# Page is set up for A4 portrait printing.
$ws.PageSetup.PaperSize = 9      # xlPaperA4
$ws.PageSetup.Orientation = 1    # xlPortrait
